# Correct linreg uptakes stderr
# - Adds "Growth-Rate-Std" (col D) and "Uptake-Std" (col F) columns to the
#   RateCompare sheet, shifting the old "Substrate-uptake" values into the
#   new column E, and fills in the corrected values for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RateCompare")
$ws.Activate() | Out-Null

# Header row
$ws.Cells.Item(1,1).Value2 = "ID"
$ws.Cells.Item(1,2).Value2 = "Tspan"
$ws.Cells.Item(1,3).Value2 = "Growth-Rate"
$ws.Cells.Item(1,4).Value2 = "Growth-Rate-Std"
$ws.Cells.Item(1,5).Value2 = "Substrate-uptake"
$ws.Cells.Item(1,6).Value2 = "Uptake-Std"

# Row 2 - 2229v1
$ws.Cells.Item(2,1).Value2 = "2229v1"
$ws.Cells.Item(2,2).Value2 = "2-8"
$ws.Cells.Item(2,3).Value2 = 0.14
$ws.Cells.Item(2,4).Value2 = 0.02
$ws.Cells.Item(2,5).Value2 = 2.99
$ws.Cells.Item(2,6).Value2 = 0.95

# Row 3 - LV3_130_v1
$ws.Cells.Item(3,1).Value2 = "LV3_130_v1"
$ws.Cells.Item(3,2).Value2 = "0-80"
$ws.Cells.Item(3,3).Value2 = 0.07
$ws.Cells.Item(3,4).Value2 = 0.01
$ws.Cells.Item(3,5).Value2 = 0.86
$ws.Cells.Item(3,6).Value2 = 0.22

# Row 4 - LV3_130_v2
$ws.Cells.Item(4,1).Value2 = "LV3_130_v2"
$ws.Cells.Item(4,2).Value2 = "0-80"
$ws.Cells.Item(4,3).Value2 = 0.07
$ws.Cells.Item(4,4).Value2 = 0.01
$ws.Cells.Item(4,5).Value2 = 0.82
$ws.Cells.Item(4,6).Value2 = 0.27

# Row 5 - LV3_200_v1
$ws.Cells.Item(5,1).Value2 = "LV3_200_v1"
$ws.Cells.Item(5,2).Value2 = "100-200"
$ws.Cells.Item(5,3).Value2 = 0.01
$ws.Cells.Item(5,4).Value2 = 0
$ws.Cells.Item(5,5).Value2 = 0.1
$ws.Cells.Item(5,6).Value2 = 0.02

# Row 6 - LV3_200_v2
$ws.Cells.Item(6,1).Value2 = "LV3_200_v2"
$ws.Cells.Item(6,2).Value2 = "0-80"
$ws.Cells.Item(6,3).Value2 = 0.06
$ws.Cells.Item(6,4).Value2 = 0.01
$ws.Cells.Item(6,5).Value2 = 1.06
$ws.Cells.Item(6,6).Value2 = 0.44

# Match the author's final selection on the sheet
$ws.Range("F7").Select() | Out-Null
